$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "Roger"
$ws.Range("C2").Value = "Comtois"

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "Sylvie"
$ws.Range("C3").Value = "Paré"

$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "Jean-Michel"
$ws.Range("C4").Value = "Esquive"

$ws.Range("A5").Value = 12
$ws.Range("B5").Value = "Thierry"
$ws.Range("C5").Value = "Plinplinplon"

$ws.Range("A6").Value = 14
$ws.Range("B6").Value = "Agzend"
$ws.Range("C6").Value = "Fireman"

$ws.Range("A7").Value = 18
$ws.Range("B7").Value = "Jean-Claude"
$ws.Range("C7").Value = "Van Damme"

$ws.Range("C8").Select()
